# Update "Stand" date on Erlauterungen sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Range("A2").Value = "Stand: 05.04.2022"


# Update data values on Symptomatische_nach_Impfstatus (sheet2)
$ws2.Range("F5").Value = 15.579489707946777
$ws2.Range("G5").Value = 4.5434503555297852
$ws2.Range("I5").Value = 3.6515777111053467
$ws2.Range("J5").Value = 0.76011365652084351
$ws2.Range("F6").Value = 23.851596832275391
$ws2.Range("G6").Value = 5.9510030746459961
$ws2.Range("I6").Value = 5.1093335151672363
$ws2.Range("J6").Value = 0.91942328214645386
$ws2.Range("F7").Value = 31.114517211914063
$ws2.Range("G7").Value = 6.8123188018798828
$ws2.Range("I7").Value = 8.1146059036254883
$ws2.Range("J7").Value = 1.698573112487793
$ws2.Range("F8").Value = 42.093276977539063
$ws2.Range("G8").Value = 8.75341796875
$ws2.Range("I8").Value = 9.5838947296142578
$ws2.Range("J8").Value = 1.9929900169372559
$ws2.Range("C9").Value = 58.885654449462891
$ws2.Range("D9").Value = 6.9515504837036133
$ws2.Range("F9").Value = 68.1337890625
$ws2.Range("G9").Value = 12.874956130981445
$ws2.Range("I9").Value = 15.466824531555176
$ws2.Range("J9").Value = 3.4228982925415039
$ws2.Range("C10").Value = 109.71903991699219
$ws2.Range("D10").Value = 12.516441345214844
$ws2.Range("F10").Value = 101.25288391113281
$ws2.Range("G10").Value = 17.367820739746094
$ws2.Range("I10").Value = 20.67054557800293
$ws2.Range("J10").Value = 4.2085738182067871
$ws2.Range("C11").Value = 139.97805786132813
$ws2.Range("D11").Value = 10.22846508026123
$ws2.Range("F11").Value = 131.39395141601563
$ws2.Range("G11").Value = 21.720512390136719
$ws2.Range("I11").Value = 31.069412231445313
$ws2.Range("J11").Value = 6.249213695526123
$ws2.Range("C12").Value = 167.66690063476563
$ws2.Range("D12").Value = 11.181144714355469
$ws2.Range("F12").Value = 136.13938903808594
$ws2.Range("G12").Value = 24.068347930908203
$ws2.Range("I12").Value = 37.572620391845703
$ws2.Range("J12").Value = 8.2231855392456055
$ws2.Range("C13").Value = 158.73031616210938
$ws2.Range("D13").Value = 11.118156433105469
$ws2.Range("F13").Value = 130.972900390625
$ws2.Range("G13").Value = 25.222780227661133
$ws2.Range("I13").Value = 42.855014801025391
$ws2.Range("J13").Value = 9.4363212585449219
$ws2.Range("C14").Value = 151.22822570800781
$ws2.Range("D14").Value = 10.951155662536621
$ws2.Range("F14").Value = 111.24652862548828
$ws2.Range("G14").Value = 22.815650939941406
$ws2.Range("I14").Value = 40.95672607421875
$ws2.Range("J14").Value = 9.5979681015014648
$ws2.Range("C15").Value = 131.73646545410156
$ws2.Range("D15").Value = 8.5513534545898438
$ws2.Range("F15").Value = 94.408668518066406
$ws2.Range("G15").Value = 21.620561599731445
$ws2.Range("I15").Value = 36.388134002685547
$ws2.Range("J15").Value = 9.6124238967895508
$ws2.Range("C16").Value = 146.70413208007813
$ws2.Range("D16").Value = 11.105525016784668
$ws2.Range("F16").Value = 98.90667724609375
$ws2.Range("G16").Value = 26.884830474853516
$ws2.Range("I16").Value = 38.986625671386719
$ws2.Range("J16").Value = 12.918940544128418
$ws2.Range("C17").Value = 147.560791015625
$ws2.Range("D17").Value = 9.4527034759521484
$ws2.Range("F17").Value = 105.81330871582031
$ws2.Range("G17").Value = 29.649946212768555
$ws2.Range("I17").Value = 48.420200347900391
$ws2.Range("J17").Value = 15.028476715087891
$ws2.Range("C18").Value = 165.58943176269531
$ws2.Range("D18").Value = 12.863425254821777
$ws2.Range("F18").Value = 120.138671875
$ws2.Range("G18").Value = 36.062320709228516
$ws2.Range("I18").Value = 54.553352355957031
$ws2.Range("J18").Value = 17.811700820922852
$ws2.Range("C19").Value = 228.16297912597656
$ws2.Range("D19").Value = 17.674701690673828
$ws2.Range("F19").Value = 177.06596374511719
$ws2.Range("G19").Value = 57.142707824707031
$ws2.Range("H19").Value = 29.9613037109375
$ws2.Range("I19").Value = 86.675483703613281
$ws2.Range("J19").Value = 31.245054244995117
$ws2.Range("K19").Value = 11.505082130432129
$ws2.Range("C20").Value = 319.3800048828125
$ws2.Range("D20").Value = 27.951986312866211
$ws2.Range("F20").Value = 228.24722290039063
$ws2.Range("G20").Value = 77.39703369140625
$ws2.Range("H20").Value = 44.699871063232422
$ws2.Range("I20").Value = 112.49837493896484
$ws2.Range("J20").Value = 46.892833709716797
$ws2.Range("K20").Value = 13.502330780029297
$ws2.Range("C21").Value = 340.99136352539063
$ws2.Range("D21").Value = 29.689472198486328
$ws2.Range("F21").Value = 268.38327026367188
$ws2.Range("G21").Value = 95.132537841796875
$ws2.Range("H21").Value = 44.969005584716797
$ws2.Range("I21").Value = 144.75904846191406
$ws2.Range("J21").Value = 59.462909698486328
$ws2.Range("K21").Value = 16.821010589599609
$ws2.Range("C22").Value = 494.7481689453125
$ws2.Range("D22").Value = 51.461414337158203
$ws2.Range("F22").Value = 329.8192138671875
$ws2.Range("G22").Value = 120.78118896484375
$ws2.Range("H22").Value = 59.339778900146484
$ws2.Range("I22").Value = 182.39974975585938
$ws2.Range("J22").Value = 76.635581970214844
$ws2.Range("K22").Value = 17.98570442199707
$ws2.Range("C23").Value = 541.16571044921875
$ws2.Range("D23").Value = 55.644096374511719
$ws2.Range("F23").Value = 359.59881591796875
$ws2.Range("G23").Value = 144.43963623046875
$ws2.Range("H23").Value = 63.652507781982422
$ws2.Range("I23").Value = 195.91624450683594
$ws2.Range("J23").Value = 88.068359375
$ws2.Range("K23").Value = 22.574935913085938
$ws2.Range("C24").Value = 548.7945556640625
$ws2.Range("D24").Value = 59.686405181884766
$ws2.Range("F24").Value = 384.30523681640625
$ws2.Range("G24").Value = 143.8603515625
$ws2.Range("H24").Value = 64.302833557128906
$ws2.Range("I24").Value = 221.16719055175781
$ws2.Range("J24").Value = 89.060745239257813
$ws2.Range("K24").Value = 25.033468246459961
$ws2.Range("C25").Value = 511.35379028320313
$ws2.Range("D25").Value = 57.741035461425781
$ws2.Range("E25").Value = 19.978689193725586
$ws2.Range("F25").Value = 349.98968505859375
$ws2.Range("G25").Value = 129.84626770019531
$ws2.Range("H25").Value = 49.186248779296875
$ws2.Range("I25").Value = 211.4039306640625
$ws2.Range("J25").Value = 75.151405334472656
$ws2.Range("K25").Value = 20.451757431030273
$ws2.Range("C26").Value = 458.69580078125
$ws2.Range("D26").Value = 55.060882568359375
$ws2.Range("E26").Value = 18.315391540527344
$ws2.Range("F26").Value = 305.02093505859375
$ws2.Range("G26").Value = 112.61207580566406
$ws2.Range("H26").Value = 39.174015045166016
$ws2.Range("I26").Value = 195.61248779296875
$ws2.Range("J26").Value = 65.246223449707031
$ws2.Range("K26").Value = 16.007625579833984
$ws2.Range("C27").Value = 378.64093017578125
$ws2.Range("D27").Value = 49.932506561279297
$ws2.Range("E27").Value = 10.906182289123535
$ws2.Range("F27").Value = 252.13395690917969
$ws2.Range("G27").Value = 96.065528869628906
$ws2.Range("H27").Value = 29.440858840942383
$ws2.Range("I27").Value = 160.36473083496094
$ws2.Range("J27").Value = 49.858909606933594
$ws2.Range("K27").Value = 11.177370071411133
$ws2.Range("C28").Value = 237.64480590820313
$ws2.Range("D28").Value = 46.775962829589844
$ws2.Range("E28").Value = 10.381396293640137
$ws2.Range("F28").Value = 170.08798217773438
$ws2.Range("G28").Value = 79.344261169433594
$ws2.Range("H28").Value = 28.509872436523438
$ws2.Range("I28").Value = 113.55418395996094
$ws2.Range("J28").Value = 36.187126159667969
$ws2.Range("K28").Value = 9.4942455291748047
$ws2.Range("C29").Value = 214.46693420410156
$ws2.Range("D29").Value = 72.118942260742188
$ws2.Range("E29").Value = 24.428625106811523
$ws2.Range("F29").Value = 175.22637939453125
$ws2.Range("G29").Value = 114.44139862060547
$ws2.Range("H29").Value = 47.393321990966797
$ws2.Range("I29").Value = 110.95841217041016
$ws2.Range("J29").Value = 42.623996734619141
$ws2.Range("K29").Value = 16.151958465576172
$ws2.Range("A30").Value = 2022
$ws2.Range("B30").Value = 1
$ws2.Range("C30").Value = 279.9549560546875
$ws2.Range("D30").Value = 153.20974731445313
$ws2.Range("E30").Value = 72.157737731933594
$ws2.Range("F30").Value = 196.48655700683594
$ws2.Range("G30").Value = 187.83407592773438
$ws2.Range("H30").Value = 92.132102966308594
$ws2.Range("I30").Value = 108.50118255615234
$ws2.Range("J30").Value = 55.630519866943359
$ws2.Range("K30").Value = 22.658170700073242
$ws2.Range("A31").Value = 2022
$ws2.Range("B31").Value = 2
$ws2.Range("C31").Value = 354.50820922851563
$ws2.Range("D31").Value = 215.03427124023438
$ws2.Range("E31").Value = 97.163429260253906
$ws2.Range("F31").Value = 187.15966796875
$ws2.Range("G31").Value = 190.16102600097656
$ws2.Range("H31").Value = 96.412445068359375
$ws2.Range("I31").Value = 77.617240905761719
$ws2.Range("J31").Value = 47.613735198974609
$ws2.Range("K31").Value = 23.064493179321289
$ws2.Range("A32").Value = 2022
$ws2.Range("B32").Value = 3
$ws2.Range("C32").Value = 450.94583129882813
$ws2.Range("D32").Value = 235.06661987304688
$ws2.Range("E32").Value = 80.451744079589844
$ws2.Range("F32").Value = 225.69975280761719
$ws2.Range("G32").Value = 209.570556640625
$ws2.Range("H32").Value = 112.92941284179688
$ws2.Range("I32").Value = 82.221992492675781
$ws2.Range("J32").Value = 53.509426116943359
$ws2.Range("K32").Value = 26.945652008056641
$ws2.Range("A33").Value = 2022
$ws2.Range("B33").Value = 4
$ws2.Range("C33").Value = 540.42767333984375
$ws2.Range("D33").Value = 255.98609924316406
$ws2.Range("E33").Value = 86.779647827148438
$ws2.Range("F33").Value = 255.21260070800781
$ws2.Range("G33").Value = 211.75428771972656
$ws2.Range("H33").Value = 118.46305084228516
$ws2.Range("I33").Value = 90.665985107421875
$ws2.Range("J33").Value = 54.088367462158203
$ws2.Range("K33").Value = 29.253576278686523
$ws2.Range("A34").Value = 2022
$ws2.Range("B34").Value = 5
$ws2.Range("C34").Value = 547.45452880859375
$ws2.Range("D34").Value = 256.13449096679688
$ws2.Range("E34").Value = 81.130302429199219
$ws2.Range("F34").Value = 269.998291015625
$ws2.Range("G34").Value = 211.25323486328125
$ws2.Range("H34").Value = 122.54711151123047
$ws2.Range("I34").Value = 111.28116607666016
$ws2.Range("J34").Value = 56.525466918945313
$ws2.Range("K34").Value = 32.386493682861328
$ws2.Range("A35").Value = 2022
$ws2.Range("B35").Value = 6
$ws2.Range("C35").Value = 490.5758056640625
$ws2.Range("D35").Value = 263.06201171875
$ws2.Range("E35").Value = 95.418014526367188
$ws2.Range("F35").Value = 259.24334716796875
$ws2.Range("G35").Value = 227.93902587890625
$ws2.Range("H35").Value = 141.06892395019531
$ws2.Range("I35").Value = 113.61431884765625
$ws2.Range("J35").Value = 65.874710083007813
$ws2.Range("K35").Value = 40.134658813476563
$ws2.Range("A36").Value = 2022
$ws2.Range("B36").Value = 7
$ws2.Range("C36").Value = 365.7757568359375
$ws2.Range("D36").Value = 232.00617980957031
$ws2.Range("E36").Value = 99.559761047363281
$ws2.Range("F36").Value = 225.4832763671875
$ws2.Range("G36").Value = 206.24375915527344
$ws2.Range("H36").Value = 140.08016967773438
$ws2.Range("I36").Value = 102.29685974121094
$ws2.Range("J36").Value = 66.625686645507813
$ws2.Range("K36").Value = 43.00048828125
$ws2.Range("A37").Value = 2022
$ws2.Range("B37").Value = 8
$ws2.Range("C37").Value = 311.01739501953125
$ws2.Range("D37").Value = 195.32931518554688
$ws2.Range("E37").Value = 108.71251678466797
$ws2.Range("F37").Value = 202.10893249511719
$ws2.Range("G37").Value = 175.56805419921875
$ws2.Range("H37").Value = 135.97560119628906
$ws2.Range("I37").Value = 97.59765625
$ws2.Range("J37").Value = 61.440284729003906
$ws2.Range("K37").Value = 44.062320709228516
$ws2.Range("A38").Value = 2022
$ws2.Range("B38").Value = 9
$ws2.Range("C38").Value = 274.38662719726563
$ws2.Range("D38").Value = 187.81086730957031
$ws2.Range("E38").Value = 133.20758056640625
$ws2.Range("F38").Value = 171.49378967285156
$ws2.Range("G38").Value = 176.67901611328125
$ws2.Range("H38").Value = 157.38577270507813
$ws2.Range("I38").Value = 87.752182006835938
$ws2.Range("J38").Value = 61.103157043457031
$ws2.Range("K38").Value = 45.558021545410156
$ws2.Range("A39").Value = 2022
$ws2.Range("B39").Value = 10
$ws2.Range("C39").Value = 273.8135986328125
$ws2.Range("D39").Value = 214.14093017578125
$ws2.Range("E39").Value = 167.59466552734375
$ws2.Range("F39").Value = 169.29322814941406
$ws2.Range("G39").Value = 181.11723327636719
$ws2.Range("H39").Value = 172.02052307128906
$ws2.Range("I39").Value = 89.202407836914063
$ws2.Range("J39").Value = 58.501285552978516
$ws2.Range("K39").Value = 51.736564636230469
$ws2.Range("A40").Value = 2022
$ws2.Range("B40").Value = 11
$ws2.Range("C40").Value = 243.9395751953125
$ws2.Range("D40").Value = 181.45491027832031
$ws2.Range("E40").Value = 134.27107238769531
$ws2.Range("F40").Value = 140.27223205566406
$ws2.Range("G40").Value = 156.82878112792969
$ws2.Range("H40").Value = 155.20498657226563
$ws2.Range("I40").Value = 81.214103698730469
$ws2.Range("J40").Value = 56.381351470947266
$ws2.Range("K40").Value = 52.724033355712891
$ws2.Range("A41").Value = 2022
$ws2.Range("B41").Value = 12
$ws2.Range("C41").Value = 202.58697509765625
$ws2.Range("D41").Value = 142.62942504882813
$ws2.Range("E41").Value = 115.64024353027344
$ws2.Range("F41").Value = 141.96604919433594
$ws2.Range("G41").Value = 122.65164184570313
$ws2.Range("H41").Value = 131.66410827636719
$ws2.Range("I41").Value = 100.51689910888672
$ws2.Range("J41").Value = 52.169082641601563
$ws2.Range("K41").Value = 47.485225677490234

# Update data values on Hospitalisierte_nach_Impfstatus (sheet3)
$ws3.Range("F5").Value = 0.8965039849281311
$ws3.Range("G5").Value = 0.11733444035053253
$ws3.Range("I5").Value = 1.1518399715423584
$ws3.Range("J5").Value = 0.27391484379768372
$ws3.Range("F6").Value = 1.3520278930664063
$ws3.Range("G6").Value = 0.10522402077913284
$ws3.Range("I6").Value = 1.7705610990524292
$ws3.Range("J6").Value = 0.32562908530235291
$ws3.Range("F7").Value = 1.8427330255508423
$ws3.Range("G7").Value = 0.18497374653816223
$ws3.Range("I7").Value = 2.9036290645599365
$ws3.Range("J7").Value = 0.43355187773704529
$ws3.Range("F8").Value = 2.5714592933654785
$ws3.Range("G8").Value = 0.26815721392631531
$ws3.Range("I8").Value = 3.7594282627105713
$ws3.Range("J8").Value = 0.49126926064491272
$ws3.Range("C9").Value = 0.91784751415252686
$ws3.Range("F9").Value = 4.5818796157836914
$ws3.Range("G9").Value = 0.3071734607219696
$ws3.Range("I9").Value = 5.5777573585510254
$ws3.Range("J9").Value = 0.89222466945648193
$ws3.Range("C10").Value = 1.5034043788909912
$ws3.Range("D10").Value = 0.21958670020103455
$ws3.Range("F10").Value = 6.6674065589904785
$ws3.Range("G10").Value = 0.47192955017089844
$ws3.Range("I10").Value = 8.1362781524658203
$ws3.Range("J10").Value = 1.1640735864639282
$ws3.Range("C11").Value = 1.6064574718475342
$ws3.Range("D11").Value = 0.1763528436422348
$ws3.Range("F11").Value = 9.5160694122314453
$ws3.Range("G11").Value = 0.53301870822906494
$ws3.Range("I11").Value = 11.486437797546387
$ws3.Range("J11").Value = 1.5933425426483154
$ws3.Range("C12").Value = 1.5599944591522217
$ws3.Range("F12").Value = 8.9320535659790039
$ws3.Range("G12").Value = 0.63688105344772339
$ws3.Range("I12").Value = 13.408603668212891
$ws3.Range("J12").Value = 2.0353915691375732
$ws3.Range("C13").Value = 1.9050359725952148
$ws3.Range("D13").Value = 0.12217754870653152
$ws3.Range("F13").Value = 8.0663919448852539
$ws3.Range("G13").Value = 0.57600516080856323
$ws3.Range("I13").Value = 15.193523406982422
$ws3.Range("J13").Value = 2.1368114948272705
$ws3.Range("C14").Value = 1.5186393260955811
$ws3.Range("D14").Value = 0.10632190108299255
$ws3.Range("F14").Value = 6.8768253326416016
$ws3.Range("G14").Value = 0.50769138336181641
$ws3.Range("I14").Value = 13.76993465423584
$ws3.Range("J14").Value = 1.9568183422088623
$ws3.Range("C15").Value = 1.023479700088501
$ws3.Range("D15").Value = 0.19003006815910339
$ws3.Range("F15").Value = 5.6586790084838867
$ws3.Range("G15").Value = 0.54533559083938599
$ws3.Range("I15").Value = 13.115262985229492
$ws3.Range("J15").Value = 2.0461592674255371
$ws3.Range("C16").Value = 1.1293619871139526
$ws3.Range("F16").Value = 5.532562255859375
$ws3.Range("G16").Value = 0.62538093328475952
$ws3.Range("I16").Value = 13.439489364624023
$ws3.Range("J16").Value = 2.4941012859344482
$ws3.Range("C17").Value = 0.88428819179534912
$ws3.Range("D17").Value = 0.07502145320177078
$ws3.Range("F17").Value = 5.2614350318908691
$ws3.Range("G17").Value = 0.57409882545471191
$ws3.Range("I17").Value = 15.192448616027832
$ws3.Range("J17").Value = 2.9793918132781982
$ws3.Range("C18").Value = 1.4882683753967285
$ws3.Range("D18").Value = 0.06699700653553009
$ws3.Range("F18").Value = 6.0481534004211426
$ws3.Range("G18").Value = 0.6710008978843689
$ws3.Range("I18").Value = 17.597856521606445
$ws3.Range("J18").Value = 3.5389840602874756
$ws3.Range("C19").Value = 1.5545696020126343
$ws3.Range("D19").Value = 0.18869787454605103
$ws3.Range("F19").Value = 7.9692573547363281
$ws3.Range("G19").Value = 1.0772727727890015
$ws3.Range("H19").Value = 1.8048977851867676
$ws3.Range("I19").Value = 28.145431518554688
$ws3.Range("J19").Value = 5.7598333358764648
$ws3.Range("K19").Value = 4.090695858001709
$ws3.Range("C20").Value = 1.9867488145828247
$ws3.Range("D20").Value = 0.23941744863986969
$ws3.Range("F20").Value = 9.448969841003418
$ws3.Range("G20").Value = 1.2750540971755981
$ws3.Range("H20").Value = 2.0318124294281006
$ws3.Range("I20").Value = 32.281051635742188
$ws3.Range("J20").Value = 7.491541862487793
$ws3.Range("K20").Value = 5.1145191192626953
$ws3.Range("C21").Value = 2.1061594486236572
$ws3.Range("D21").Value = 0.22970578074455261
$ws3.Range("F21").Value = 9.9404172897338867
$ws3.Range("G21").Value = 1.3846265077590942
$ws3.Range("H21").Value = 2.0754926204681396
$ws3.Range("I21").Value = 38.53094482421875
$ws3.Range("J21").Value = 8.526580810546875
$ws3.Range("K21").Value = 5.5517621040344238
$ws3.Range("C22").Value = 2.7016963958740234
$ws3.Range("D22").Value = 0.38776093721389771
$ws3.Range("F22").Value = 12.022501945495605
$ws3.Range("G22").Value = 1.681121826171875
$ws3.Range("H22").Value = 0.53459256887435913
$ws3.Range("I22").Value = 48.531871795654297
$ws3.Range("J22").Value = 10.165729522705078
$ws3.Range("K22").Value = 6.0168004035949707
$ws3.Range("C23").Value = 3.686119556427002
$ws3.Range("D23").Value = 0.75049841403961182
$ws3.Range("F23").Value = 13.322453498840332
$ws3.Range("G23").Value = 2.0465779304504395
$ws3.Range("H23").Value = 2.0368802547454834
$ws3.Range("I23").Value = 55.687202453613281
$ws3.Range("J23").Value = 10.863153457641602
$ws3.Range("K23").Value = 7.3968939781188965
$ws3.Range("C24").Value = 4.1232419013977051
$ws3.Range("D24").Value = 0.41775259375572205
$ws3.Range("F24").Value = 14.197780609130859
$ws3.Range("G24").Value = 2.0450842380523682
$ws3.Range("H24").Value = 1.4877262115478516
$ws3.Range("I24").Value = 63.325370788574219
$ws3.Range("J24").Value = 11.107912063598633
$ws3.Range("K24").Value = 6.4143924713134766
$ws3.Range("C25").Value = 3.3262591361999512
$ws3.Range("D25").Value = 0.46069973707199097
$ws3.Range("E25").Value = 0
$ws3.Range("F25").Value = 13.703057289123535
$ws3.Range("G25").Value = 1.9836174249649048
$ws3.Range("H25").Value = 1.6283664703369141
$ws3.Range("I25").Value = 62.060493469238281
$ws3.Range("J25").Value = 9.8186845779418945
$ws3.Range("K25").Value = 5.1852922439575195
$ws3.Range("C26").Value = 3.042595386505127
$ws3.Range("D26").Value = 0.4542144238948822
$ws3.Range("E26").Value = 0
$ws3.Range("F26").Value = 12.218108177185059
$ws3.Range("G26").Value = 2.0909490585327148
$ws3.Range("H26").Value = 0.95357799530029297
$ws3.Range("I26").Value = 64.379653930664063
$ws3.Range("J26").Value = 9.7886390686035156
$ws3.Range("K26").Value = 4.8287463188171387
$ws3.Range("C27").Value = 3.1771197319030762
$ws3.Range("D27").Value = 0.35022801160812378
$ws3.Range("E27").Value = 0
$ws3.Range("F27").Value = 10.32520580291748
$ws3.Range("G27").Value = 1.59928297996521
$ws3.Range("H27").Value = 0.79651123285293579
$ws3.Range("I27").Value = 52.8101806640625
$ws3.Range("J27").Value = 8.0745458602905273
$ws3.Range("K27").Value = 3.251598596572876
$ws3.Range("C28").Value = 1.9835004806518555
$ws3.Range("D28").Value = 0.49708780646324158
$ws3.Range("E28").Value = 0
$ws3.Range("F28").Value = 7.6170430183410645
$ws3.Range("G28").Value = 1.3434337377548218
$ws3.Range("H28").Value = 0.46413716673851013
$ws3.Range("I28").Value = 38.460765838623047
$ws3.Range("J28").Value = 5.1536784172058105
$ws3.Range("K28").Value = 1.6722300052642822
$ws3.Range("C29").Value = 1.4821109771728516
$ws3.Range("D29").Value = 0.74044090509414673
$ws3.Range("E29").Value = 0
$ws3.Range("F29").Value = 7.6601696014404297
$ws3.Range("G29").Value = 1.7128785848617554
$ws3.Range("H29").Value = 0.58790099620819092
$ws3.Range("I29").Value = 37.600406646728516
$ws3.Range("J29").Value = 5.6244697570800781
$ws3.Range("K29").Value = 1.7773818969726563
$ws3.Range("A30").Value = 2022
$ws3.Range("B30").Value = 1
$ws3.Range("C30").Value = 2.2342774868011475
$ws3.Range("D30").Value = 1.065638542175293
$ws3.Range("E30").Value = 0
$ws3.Range("F30").Value = 6.3414874076843262
$ws3.Range("G30").Value = 2.2682125568389893
$ws3.Range("H30").Value = 0.79344689846038818
$ws3.Range("I30").Value = 29.332340240478516
$ws3.Range("J30").Value = 5.5113825798034668
$ws3.Range("K30").Value = 1.5924555063247681
$ws3.Range("A31").Value = 2022
$ws3.Range("B31").Value = 2
$ws3.Range("C31").Value = 2.9203550815582275
$ws3.Range("D31").Value = 1.2674882411956787
$ws3.Range("E31").Value = 0.67241132259368896
$ws3.Range("F31").Value = 4.9700899124145508
$ws3.Range("G31").Value = 2.2266185283660889
$ws3.Range("H31").Value = 0.94868725538253784
$ws3.Range("I31").Value = 20.174749374389648
$ws3.Range("J31").Value = 4.0606174468994141
$ws3.Range("K31").Value = 1.7148321866989136
$ws3.Range("A32").Value = 2022
$ws3.Range("B32").Value = 3
$ws3.Range("C32").Value = 2.7625722885131836
$ws3.Range("D32").Value = 1.5018576383590698
$ws3.Range("E32").Value = 0.21626812219619751
$ws3.Range("F32").Value = 5.5170822143554688
$ws3.Range("G32").Value = 2.2224414348602295
$ws3.Range("H32").Value = 0.92679446935653687
$ws3.Range("I32").Value = 16.285514831542969
$ws3.Range("J32").Value = 5.1789765357971191
$ws3.Range("K32").Value = 2.0401084423065186
$ws3.Range("A33").Value = 2022
$ws3.Range("B33").Value = 4
$ws3.Range("C33").Value = 3.6706478595733643
$ws3.Range("D33").Value = 1.1498841047286987
$ws3.Range("E33").Value = 0.61545848846435547
$ws3.Range("F33").Value = 5.4210305213928223
$ws3.Range("G33").Value = 2.0890429019927979
$ws3.Range("H33").Value = 0.7637590765953064
$ws3.Range("I33").Value = 20.503454208374023
$ws3.Range("J33").Value = 5.8318877220153809
$ws3.Range("K33").Value = 1.9167156219482422
$ws3.Range("A34").Value = 2022
$ws3.Range("B34").Value = 5
$ws3.Range("C34").Value = 3.737518310546875
$ws3.Range("D34").Value = 0.98938566446304321
$ws3.Range("E34").Value = 0.35687816143035889
$ws3.Range("F34").Value = 5.805516242980957
$ws3.Range("G34").Value = 1.8839800357818604
$ws3.Range("H34").Value = 0.81595903635025024
$ws3.Range("I34").Value = 24.302783966064453
$ws3.Range("J34").Value = 7.2750115394592285
$ws3.Range("K34").Value = 2.3932876586914063
$ws3.Range("A35").Value = 2022
$ws3.Range("B35").Value = 6
$ws3.Range("C35").Value = 4.0932998657226563
$ws3.Range("D35").Value = 1.3518719673156738
$ws3.Range("E35").Value = 0.40517202019691467
$ws3.Range("F35").Value = 5.2111330032348633
$ws3.Range("G35").Value = 1.9493025541305542
$ws3.Range("H35").Value = 0.82344377040863037
$ws3.Range("I35").Value = 27.358386993408203
$ws3.Range("J35").Value = 7.2383108139038086
$ws3.Range("K35").Value = 2.5349380970001221
$ws3.Range("A36").Value = 2022
$ws3.Range("B36").Value = 7
$ws3.Range("C36").Value = 1.9396251440048218
$ws3.Range("D36").Value = 0.76234233379364014
$ws3.Range("E36").Value = 0.3683987557888031
$ws3.Range("F36").Value = 4.5480761528015137
$ws3.Range("G36").Value = 1.5820398330688477
$ws3.Range("H36").Value = 0.78709602355957031
$ws3.Range("I36").Value = 25.886768341064453
$ws3.Range("J36").Value = 6.6495299339294434
$ws3.Range("K36").Value = 2.801030158996582
$ws3.Range("A37").Value = 2022
$ws3.Range("B37").Value = 8
$ws3.Range("C37").Value = 2.3299543857574463
$ws3.Range("D37").Value = 0.58268606662750244
$ws3.Range("E37").Value = 0.08669260144233704
$ws3.Range("F37").Value = 3.8013887405395508
$ws3.Range("G37").Value = 1.5181224346160889
$ws3.Range("H37").Value = 0.84484755992889404
$ws3.Range("I37").Value = 22.593414306640625
$ws3.Range("J37").Value = 6.6647090911865234
$ws3.Range("K37").Value = 2.9476399421691895
$ws3.Range("A38").Value = 2022
$ws3.Range("B38").Value = 9
$ws3.Range("C38").Value = 2.657804012298584
$ws3.Range("D38").Value = 0.98090630769729614
$ws3.Range("E38").Value = 0.41549462080001831
$ws3.Range("F38").Value = 3.2165038585662842
$ws3.Range("G38").Value = 1.681193470954895
$ws3.Range("H38").Value = 0.92393654584884644
$ws3.Range("I38").Value = 21.088592529296875
$ws3.Range("J38").Value = 5.752161979675293
$ws3.Range("K38").Value = 2.8332104682922363
$ws3.Range("A39").Value = 2022
$ws3.Range("B39").Value = 10
$ws3.Range("C39").Value = 1.8415106534957886
$ws3.Range("D39").Value = 0.78921717405319214
$ws3.Range("E39").Value = 0.40306556224822998
$ws3.Range("F39").Value = 2.9885284900665283
$ws3.Range("G39").Value = 1.5680297613143921
$ws3.Range("H39").Value = 0.78931856155395508
$ws3.Range("I39").Value = 19.34197998046875
$ws3.Range("J39").Value = 5.6251239776611328
$ws3.Range("K39").Value = 2.9597456455230713
$ws3.Range("A40").Value = 2022
$ws3.Range("B40").Value = 11
$ws3.Range("C40").Value = 2.4209203720092773
$ws3.Range("D40").Value = 0.5938524603843689
$ws3.Range("E40").Value = 0.31556069850921631
$ws3.Range("F40").Value = 2.9298849105834961
$ws3.Range("G40").Value = 1.1344213485717773
$ws3.Range("H40").Value = 0.74979996681213379
$ws3.Range("I40").Value = 19.664699554443359
$ws3.Range("J40").Value = 7.1730036735534668
$ws3.Range("K40").Value = 2.742375373840332
$ws3.Range("A41").Value = 2022
$ws3.Range("B41").Value = 12
$ws3.Range("C41").Value = 1.4059664011001587
$ws3.Range("D41").Value = 0.33138808608055115
$ws3.Range("E41").Value = 0.30981981754302979
$ws3.Range("F41").Value = 2.3262627124786377
$ws3.Range("G41").Value = 1.0401909351348877
$ws3.Range("H41").Value = 0.59026157855987549
$ws3.Range("I41").Value = 17.685932159423828
$ws3.Range("J41").Value = 5.3557615280151367
$ws3.Range("K41").Value = 2.1509358882904053

# Update sheet views: zoom, selection, active tab
$ws3.Activate()
$excel.ActiveWindow.Zoom = 90
$ws3.Range("A3").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 90
$ws2.Range("A3").Select()
